$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("file")

$ws.Range("A6").Value = "convertNegLog10Pvalue"
$ws.Range("C6").Value = "set to TRUE if your p_value values need converting from -log10(pvalue)"

$ws.Activate()
$ws.Range("I6").Select()
